# edit.ps1 -- applies the "mcast leave stage / timeout reasons" rewrite
# described in the commit diff to the Preparation-stage bullet list.
#
# Summary of the textual change:
#  1. The paragraph ending in
#       "...If it stops due to previous reason it aborts, otherwise it
#        does the following:"
#     becomes
#       "...The mcast may leave this stage due to any of these three
#        reasons:"
#     and the (hidden) "_GoBack" bookmark is relocated into the middle
#     of the word "these" (after "the", before "se").
#  2. The sub-bullet that used to read
#       "Machine 0: It maintains a vector table of every other
#        machines' rcved field in their packet received. If all of
#        them are 1, go to transmission stage."
#     is replaced with a brand-new sentence:
#       "If it stops due to timeout, it aborts. Otherwise we have
#        situation ii and iii."
#  3. A new sub-bullet (same list level/style) is inserted right after
#     it, carrying the *old* "Machine 0: ..." text verbatim (with
#     "rcved" italicised), placed just before the "Others:" bullet.
#  4. The "Others:" bullet gains a trailing clause:
#       "Others: It goes to Transmission stage." ->
#       "Others: It goes to Transmission stage when it receives a
#        TOKEN from last machine."

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: duplicate the "Machine 0: ..." bullet into a new paragraph
# that will sit right after it (right before "Others: ..."), *before*
# we touch its original text, since we want a verbatim copy.
# ---------------------------------------------------------------------

$machineIdx = -1
$paras = @($d.Paragraphs)
for ($i = 0; $i -lt $paras.Count; $i++) {
    if ($paras[$i].Range.Text -like "Machine 0: It maintains*") {
        $machineIdx = $i
        break
    }
}

if ($machineIdx -ge 0) {
    $srcPara = $paras[$machineIdx]
    $srcPara.Range.InsertParagraphAfter() | Out-Null

    # Re-fetch paragraphs; the newly inserted paragraph is empty and
    # sits immediately after the original "Machine 0: ..." paragraph.
    $paras2 = @($d.Paragraphs)
    $newPara = $paras2[$machineIdx + 1]
    $newPara.Range.Text = "Machine 0: It maintains a vector table of every other machines’ rcved field in their packet received. If all of them are 1, go to transmission stage."

    # Re-italicise "rcved" inside the freshly written copy.
    $paras3 = @($d.Paragraphs)
    $copyPara = $paras3[$machineIdx + 1]
    $copyRng = $copyPara.Range
    if ($copyRng.Find.Execute("rcved")) {
        $copyRng.Italic = 1
    }

    # Replace the text of the *original* paragraph with the new
    # "If it stops due to timeout..." sentence (keep the paragraph
    # mark itself intact by stopping one character short of End).
    $paras4 = @($d.Paragraphs)
    $origPara = $paras4[$machineIdx]
    $origRange = $origPara.Range
    $bodyRng = $d.Range($origRange.Start, $origRange.End - 1)
    $bodyRng.Text = "If it stops due to timeout, it aborts. Otherwise we have situation ii and iii."
}

# ---------------------------------------------------------------------
# Step 2: rewrite the lead-in sentence and relocate the "_GoBack"
# bookmark into the middle of "these".
# ---------------------------------------------------------------------

$d.Content.Find.Execute(
    "If it stops due to previous reason it aborts, otherwise it does the following:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The mcast may leave this stage due to any of these three reasons:",
    2) | Out-Null

$locate = $d.Content
if ($locate.Find.Execute("any of the")) {
    $bmPos = $locate.End
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

# ---------------------------------------------------------------------
# Step 3: extend the "Others: ..." bullet with the missing clause.
# ---------------------------------------------------------------------

$othersRng = $d.Content
if ($othersRng.Find.Execute("Others: It goes to Transmission stage.")) {
    $tail = $d.Range($othersRng.End - 1, $othersRng.End)
    $tail.Text = " when it receives a TOKEN from last machine."
}
